# Updated batch file processing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 35
$ws.Range("B3").Value = 42
$ws.Range("B4").Value = 28
$ws.Range("B5").Value = 15
$ws.Range("B6").Value = 50
